# Implements: SheetCellName, SheetCellRange, SheetColumnRange
#
# Adds a new worksheet "named_ranges" at the end of the workbook with a
# single-cell value at A2 and a 2x2 range of values at B4:C5, then defines
# two workbook-level named ranges pointing at them:
#   - SINGLE_CELL -> named_ranges!$A$2
#   - RANGE_B4C5  -> named_ranges!$B$4:$C$5

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so it lands at the end
# of the tab order (Worksheets.Add() with no args inserts at the front).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "named_ranges"

# Single cell value.
$ws.Range("A2").Value = "single cell A2"

# 2x2 range, every cell sharing the same text value.
$ws.Range("B4:C5").Value = "range B4:C5"

# Widen column A so the "single cell A2" label is fully visible.
$ws.Columns.Item(1).ColumnWidth = 13.830729166666666

# Workbook-level defined names referencing the new sheet.
$wb.Names.Add("RANGE_B4C5", "=named_ranges!`$B`$4:`$C`$5")
$wb.Names.Add("SINGLE_CELL", "=named_ranges!`$A`$2")
